# Fruta / hortaliza, semanal
# Insert this week's two new Kiwi price rows (Vega Monumental Concepción)
# right before the existing row 80 block, shifting the remaining data down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("80:81").Insert()

# New row 80 - Primera
$ws.Range("A80").Value() = 11
$ws.Range("B80").Value() = "Vega Monumental Concepción"
$ws.Range("C80").Value() = "Bíobío"
$ws.Range("D80").Value() = 44449
$ws.Range("E80").Value() = 8
$ws.Range("F80").Value() = "Fruta"
$ws.Range("G80").Value() = 100101
$ws.Range("H80").Value() = "Berries"
$ws.Range("I80").Value() = 100101007
$ws.Range("J80").Value() = "Kiwi"
$ws.Range("K80").Value() = "Hayward"
$ws.Range("L80").Value() = "Primera"
$ws.Range("M80").Value() = 100
$ws.Range("N80").Value() = 12000
$ws.Range("O80").Value() = 13000
$ws.Range("P80").Value() = 12500
$ws.Range("Q80").Value() = "$/bandeja 18 kilos"
$ws.Range("R80").Value() = "Región de O'Higgins"
$ws.Range("S80").Value() = 694
$ws.Range("T80").Value() = 18

# New row 81 - Segunda
$ws.Range("A81").Value() = 11
$ws.Range("B81").Value() = "Vega Monumental Concepción"
$ws.Range("C81").Value() = "Bíobío"
$ws.Range("D81").Value() = 44449
$ws.Range("E81").Value() = 8
$ws.Range("F81").Value() = "Fruta"
$ws.Range("G81").Value() = 100101
$ws.Range("H81").Value() = "Berries"
$ws.Range("I81").Value() = 100101007
$ws.Range("J81").Value() = "Kiwi"
$ws.Range("K81").Value() = "Hayward"
$ws.Range("L81").Value() = "Segunda"
$ws.Range("M81").Value() = 50
$ws.Range("N81").Value() = 10000
$ws.Range("O81").Value() = 10000
$ws.Range("P81").Value() = 10000
$ws.Range("Q81").Value() = "$/bandeja 18 kilos"
$ws.Range("R81").Value() = "Región de O'Higgins"
$ws.Range("S81").Value() = 556
$ws.Range("T81").Value() = 18
